# Update Name of Algo
# Applies corrected imputed values to columns C and D on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Cell = "D11";  Value = -7.139 },
    @{ Cell = "C12";  Value = -10.534 },
    @{ Cell = "D23";  Value = -7.927 },
    @{ Cell = "D28";  Value = -7.834000000000001 },
    @{ Cell = "C32";  Value = -11.6 },
    @{ Cell = "D32";  Value = -7.540999999999999 },
    @{ Cell = "D34";  Value = -7.933 },
    @{ Cell = "C36";  Value = -12.732 },
    @{ Cell = "C38";  Value = -12.444 },
    @{ Cell = "D42";  Value = -8.347 },
    @{ Cell = "C46";  Value = -13.935 },
    @{ Cell = "C54";  Value = -12.429 },
    @{ Cell = "D54";  Value = -7.531000000000001 },
    @{ Cell = "C55";  Value = -13.65 },
    @{ Cell = "C67";  Value = -11.651 },
    @{ Cell = "C69";  Value = -10.927 },
    @{ Cell = "C72";  Value = -11.753 },
    @{ Cell = "C91";  Value = -12.173 },
    @{ Cell = "D97";  Value = -7.674000000000001 },
    @{ Cell = "C99";  Value = -11.242 },
    @{ Cell = "D99";  Value = -7.428 },
    @{ Cell = "D101"; Value = -7.656999999999999 },
    @{ Cell = "C104"; Value = -12.381 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
